$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write new row 79 data (new day's values)
$ws.Cells.Item(79, 1).Value = 45819
$ws.Cells.Item(79, 2).Value = 336
$ws.Cells.Item(79, 3).Value = 338
$ws.Cells.Item(79, 4).Value = 343

# The "latest row" date cell (column A) carries a distinct date-only
# number format; it moves from the previous last row (78) to the new
# last row (79). Row 78's A cell reverts to the standard datetime format
# used by all the other (non-latest) rows.
$ws.Cells.Item(78, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79, 1).NumberFormat = "YYYY-MM-DD"
